# update for 7 march 2022
#
# The "current week" reading in the mobility tracker (Sheet1 row 8) moves
# from 23-Feb-22 to 03-Feb-22, and the Transit Mobility / Grocery Mobility
# readings for that row are revised. The embedded chart reads its series
# directly from this row, so refreshing it here keeps everything in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the underlying data table -------------------------------------
$ws.Range("A8").Value = 44595   # 03-Feb-22 (was 23-Feb-22 / 44615)
$ws.Range("C8").Value = 0.19    # Transit Mobility (was 0.17)
$ws.Range("E8").Value = 0.5     # Grocery Mobility (was 0.48)

# --- Make sure the embedded chart (sourced from this row) is current ------
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh() | Out-Null
$wb.RefreshAll() | Out-Null
$excel.CalculateFull() | Out-Null

# --- Match the author's on-save cursor / scroll position -------------------
$ws.Range("A1:F20").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("F8").Select() | Out-Null
